$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-02-02 Monday" "2026-02-03 Tuesday"

Replace-Text "717×5=" "224×5="
Replace-Text "210×7=" "823×9="
Replace-Text "368×4=" "389×5="
Replace-Text "691×7=" "548×6="
Replace-Text "429×9=" "512×4="

Replace-Text "701×3=" "892×4="
Replace-Text "310×4=" "605×9="
Replace-Text "243×9=" "251×2="
Replace-Text "298×9=" "523×6="
Replace-Text "408×2=" "930×5="

Replace-Text "617×7=" "591×4="
Replace-Text "399×9=" "662×9="
Replace-Text "713×5=" "536×6="
Replace-Text "398×4=" "156×7="
Replace-Text "194×5=" "123×4="

Replace-Text "836×6=" "303×6="
Replace-Text "113×9=" "396×6="
Replace-Text "133×2=" "927×9="
Replace-Text "955×5=" "221×8="
Replace-Text "880×3=" "145×5="

Replace-Text "110×8=" "928×7="
Replace-Text "779×8=" "549×7="
Replace-Text "487×8=" "707×6="
Replace-Text "361×5=" "142×8="
Replace-Text "638×7=" "708×4="
